$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Octubre de 2020 a las 07:37"

# --- Israel: refreshed case counts (no row move) ---
$ws.Range("B27").Value = 282872
$ws.Range("C27").Value = 1391
$ws.Range("D27").Value = 219998
$ws.Range("E27").Value = 61050

# --- Uzbekistan overtakes Nigeria: rows 58/59 swap country + data ---
$ws.Range("A58").Value = "Uzbekistan"
$ws.Range("B58").Value = 59905
$ws.Range("C58").Value = 326
$ws.Range("D58").Value = 56568
$ws.Range("E58").Value = 2843
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 494

$ws.Range("A59").Value = "Nigeria"
$ws.Range("B59").Value = 59738
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 51403
$ws.Range("E59").Value = 7222
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 1113

# --- Libano: refreshed case counts (no row move) ---
$ws.Range("B67").Value = 48342
$ws.Range("C67").Value = 245
$ws.Range("D67").Value = 43957
$ws.Range("E67").Value = 3312
$ws.Range("G67").Value = 4
$ws.Range("H67").Value = 1073

# --- Nueva Caledonia overtakes Santa Lucia: rows 207/208 swap country (data identical) ---
$ws.Range("A207").Value = "Nueva Caledonia"
$ws.Range("A208").Value = "Santa Lucia"

# --- Montserrat overtakes Islas Malvinas: rows 215/216 swap country + data ---
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
